$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values for the HORTA district result row according to the new totals.
$ws.Range("H2").Value  = 163
$ws.Range("I2").Value  = 360
$ws.Range("J2").Value  = 1627
$ws.Range("K2").Value  = 6
$ws.Range("L2").Value  = 455
$ws.Range("M2").Value  = 31
$ws.Range("O2").Value  = 0
$ws.Range("P2").Value  = 10
$ws.Range("Q2").Value  = 3
$ws.Range("R2").Value  = 19
$ws.Range("S2").Value  = 195
$ws.Range("T2").Value  = 284
$ws.Range("U2").Value  = 23
$ws.Range("V2").Value  = 2421
$ws.Range("W2").Value  = 0
$ws.Range("X2").Value  = 2496
$ws.Range("Y2").Value  = 1
$ws.Range("Z2").Value  = 34
